# The XLSForm's "survey" sheet contained an unsupported "geopoint" question
# (name "survey_gps") on row 10. Per the commit message, the example form is
# updated to no longer use this unsupported input type, so the whole survey
# row is removed and everything below it shifts up by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the entire "geopoint" row (type/name/label = geopoint/survey_gps/
# "4.Please collect the GPS co-ordinate "). Deleting the row shifts every
# subsequent row up by one and Excel keeps all formula/shared-string
# bookkeeping in sync automatically.
$ws.Rows.Item(10).Delete()

# The conditional formatting that used to sit on I18 (the "calculate"/
# "position" helper row) now lives on I17 after the shift; keep it anchored
# to the correct cell.
$ws.Range("I18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I17"))

# Restore the (shifted) active selection on the sheet.
$ws.Range("F21").Select()
